$d = $word.ActiveDocument

# --- 1. Text-only replacements (safe first; do not change paragraph count) ---

$d.Content.Find.Execute(
    "As a user, I want to be able to see all services offered by the center.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "As a visitor to the site, I want to be able to see all services offered by the center.",
    2) | Out-Null

$d.Content.Find.Execute(
    "As a user, I want to be able to see the location of the center using a Google maps API.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "As a visitor to the site, I want to be able to see the location of the center using a Google maps API.",
    2) | Out-Null

$d.Content.Find.Execute(
    "As a user, I want to be able to view and download the 100 questions for passing the citizen exam.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "As a visitor to the site, I want to be able to view the 100 questions for passing the citizen exam.",
    2) | Out-Null

$d.Content.Find.Execute(
    "As a user, I want to be able to take a practice test and see results.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "As a visitor to the site, I want to be able to take a practice test and see results (10 questions).",
    2) | Out-Null

$d.Content.Find.Execute(
    "As a user, I want to be able to view English language videos through the YouTube API.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "As a visitor to the site, I want to be able to view English language videos through the YouTube API.",
    2) | Out-Null

$d.Content.Find.Execute(
    "As an admin, I want to be notified if someone requests information about the center.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "As an admin, I want to be able to receive an email if someone requests information about the center. (built into the asp.net contact page)",
    2) | Out-Null

$d.Content.Find.Execute(
    "As an admin, I want to be able to view up coming events.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "As an admin, I want to be able to view and post up coming events.",
    2) | Out-Null

# --- 2. Insert a new bullet paragraph after the "YouTube API" item ---

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("As a visitor to the site, I want to be able to view English language videos through the YouTube API.")) {
        $p.Range.InsertParagraphAfter()
        break
    }
}

# Find the freshly-inserted empty paragraph (immediately after the YouTube one) and set its text.
$found = $false
foreach ($p in $d.Paragraphs) {
    if ($found) {
        $p.Range.Text = "As a visitor to the site, I want to be able to view upcoming events."
        $found = $false
        break
    }
    if ($p.Range.Text.StartsWith("As a visitor to the site, I want to be able to view English language videos through the YouTube API.")) {
        $found = $true
    }
}

# --- 3. Delete the "notified if someone takes a practice test" paragraph entirely ---

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("As an admin, I want to be notified if someone takes a practice test.")) {
        $p.Range.Delete()
        break
    }
}

# --- 4. Delete the first of the two trailing empty paragraphs (right after the "up coming events" bullet) ---

$found = $false
foreach ($p in $d.Paragraphs) {
    if ($found) {
        $p.Range.Delete()
        $found = $false
        break
    }
    if ($p.Range.Text.StartsWith("As an admin, I want to be able to view and post up coming events.")) {
        $found = $true
    }
}
